$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 33 (diff @ 2243)
$ws.Range("H33").Value = 847.375
$ws.Range("I33").Value = 840.3684
$ws.Range("J33").Value = 874
$ws.Range("K33").Value = 840.3684
$ws.Range("L33").Value = 874
$ws.Range("M33").Value = -611.3684
$ws.Range("N33").Value = -1332
# row 96 (diff @ 5384)
$ws.Range("H96").Value = 7143980
$ws.Range("I96").Value = 8929225
$ws.Range("J96").Value = 2998.5
$ws.Range("K96").Value = 26787675
$ws.Range("L96").Value = 8995.5
$ws.Range("M96").Value = -26786302
$ws.Range("N96").Value = -11741.5
# row 106 (diff @ 5892)
$ws.Range("H106").Value = 2857
$ws.Range("J106").Value = 3750
$ws.Range("L106").Value = 3750
$ws.Range("N106").Value = -5012
# row 113 (diff @ 6241)
$ws.Range("H113").Value = 8412.833000000001
$ws.Range("I113").Value = 4005
$ws.Range("J113").Value = 11939.1
$ws.Range("K113").Value = 4005
$ws.Range("L113").Value = 11939.1
$ws.Range("M113").Value = -751
$ws.Range("N113").Value = -18447.1
# row 114 (diff @ 6293)
$ws.Range("H114").Value = 58999.5
$ws.Range("J114").Value = 58999.5
$ws.Range("L114").Value = 58999.5
$ws.Range("N114").Value = -67677.5
# row 138 (diff @ 7493)
$ws.Range("H138").Value = 3796.5652
$ws.Range("J138").Value = 4151.5
$ws.Range("L138").Value = 12454.5
$ws.Range("N138").Value = -22734.5

$ws = $wb.Worksheets.Item("ARM")
# row 32 (diff @ 9271)
$ws.Range("H32").Value = 26304.857
$ws.Range("I32").Value = 27975.35
$ws.Range("K32").Value = 27975.35
$ws.Range("M32").Value = -27688.35
# row 37 (diff @ 9525)
$ws.Range("H37").Value = 22856.715
$ws.Range("J37").Value = 33332.332
$ws.Range("L37").Value = 33332.332
$ws.Range("N37").Value = -33878.332
# row 43 (diff @ 9819)
$ws.Range("H43").Value = 69914.664
$ws.Range("J43").Value = 92372
$ws.Range("L43").Value = 92372
$ws.Range("N43").Value = -92998
# row 45 (diff @ 9920)
$ws.Range("H45").Value = 2884.6924
$ws.Range("I45").Value = 2092.353
$ws.Range("K45").Value = 2092.353
$ws.Range("M45").Value = -1715.353
# row 55 (diff @ 10407)
$ws.Range("H55").Value = 42998.668
$ws.Range("J55").Value = 42998.668
$ws.Range("L55").Value = 42998.668
$ws.Range("N55").Value = -43628.668
# row 61 (diff @ 10686)
$ws.Range("H61").Value = 2690.2258
$ws.Range("I61").Value = 1292.2307
$ws.Range("K61").Value = 1292.2307
$ws.Range("M61").Value = -1080.2307
# row 74 (diff @ 11308)
$ws.Range("H74").Value = 463498.53
$ws.Range("I74").Value = 546634.6
$ws.Range("K74").Value = 546634.6
$ws.Range("M74").Value = -545760.6
# row 77 (diff @ 11458)
$ws.Range("H77").Value = 463498.53
$ws.Range("I77").Value = 546634.6
$ws.Range("K77").Value = 2733173
$ws.Range("M77").Value = -2728805
# row 132 (diff @ 14153)
$ws.Range("H132").Value = 1335.0714
$ws.Range("I132").Value = 1053.2307
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 3159.6921
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -629.6921000000002
$ws.Range("N132").Value = -20057
# row 136 (diff @ 14346)
$ws.Range("H136").Value = 2690.2258
$ws.Range("I136").Value = 1292.2307
$ws.Range("K136").Value = 3876.6921
$ws.Range("M136").Value = -1326.6921

$ws = $wb.Worksheets.Item("BSM")
# row 99 (diff @ 19481)
$ws.Range("H99").Value = 2033.7059
$ws.Range("I99").Value = 1528.0769
$ws.Range("K99").Value = 1528.0769
$ws.Range("M99").Value = -30.07690000000002

$ws = $wb.Worksheets.Item("CRP")
# row 16 (diff @ 22350)
$ws.Range("H16").Value = 1600.5
$ws.Range("I16").Value = 1637.5
$ws.Range("J16").Value = 1489.5
$ws.Range("K16").Value = 1637.5
$ws.Range("L16").Value = 1489.5
$ws.Range("M16").Value = -1350.5
$ws.Range("N16").Value = -2063.5
# row 31 (diff @ 23103)
$ws.Range("H31").Value = 5884290
$ws.Range("I31").Value = 6668195.5
$ws.Range("K31").Value = 6668195.5
$ws.Range("M31").Value = -6667900.5
# row 34 (diff @ 23253)
$ws.Range("H34").Value = 5884290
$ws.Range("I34").Value = 6668195.5
$ws.Range("K34").Value = 6668195.5
$ws.Range("M34").Value = -6667993.5
# row 94 (diff @ 26223)
$ws.Range("H94").Value = 1337.4
$ws.Range("J94").Value = 1399.5
$ws.Range("L94").Value = 1399.5
$ws.Range("N94").Value = -2301.5
# row 113 (diff @ 27160)
$ws.Range("H113").Value = 1600.5
$ws.Range("I113").Value = 1637.5
$ws.Range("J113").Value = 1489.5
$ws.Range("K113").Value = 1637.5
$ws.Range("L113").Value = 1489.5
$ws.Range("M113").Value = 532.5
$ws.Range("N113").Value = -5829.5
# row 132 (diff @ 28079)
$ws.Range("H132").Value = 28230.682
$ws.Range("I132").Value = 30964.9
$ws.Range("J132").Value = 888.5
$ws.Range("K132").Value = 92894.70000000001
$ws.Range("L132").Value = 2665.5
$ws.Range("M132").Value = -90364.70000000001
$ws.Range("N132").Value = -7725.5

$ws = $wb.Worksheets.Item("CUL")
# row 37 (diff @ 30429)
$ws.Range("H37").Value = 39934.39
$ws.Range("J37").Value = 39934.39
$ws.Range("L37").Value = 119803.17
$ws.Range("N37").Value = -120027.17
# row 121 (diff @ 34671)
$ws.Range("H121").Value = 151361.28
$ws.Range("I121").Value = 209893.6
$ws.Range("J121").Value = 5030.5
$ws.Range("K121").Value = 629680.8
$ws.Range("L121").Value = 15091.5
$ws.Range("M121").Value = -628370.8
$ws.Range("N121").Value = -17711.5
# row 122 (diff @ 34723)
$ws.Range("H122").Value = 89396.914
$ws.Range("J122").Value = 1241.7742
$ws.Range("L122").Value = 11175.9678
$ws.Range("N122").Value = -16075.9678
# row 131 (diff @ 35185)
$ws.Range("H131").Value = 165346.58
$ws.Range("J131").Value = 1968.0526
$ws.Range("L131").Value = 5904.1578
$ws.Range("N131").Value = -15984.1578
# row 133 (diff @ 35289)
$ws.Range("H133").Value = 6177
$ws.Range("J133").Value = 5000
$ws.Range("L133").Value = 15000
$ws.Range("N133").Value = -25120

$ws = $wb.Worksheets.Item("GSM")
# row 80 (diff @ 39652)
$ws.Range("H80").Value = 6276.9414
$ws.Range("I80").Value = 3050.1667
$ws.Range("K80").Value = 3050.1667
$ws.Range("M80").Value = -2052.1667
# row 83 (diff @ 39799)
$ws.Range("H83").Value = 6276.9414
$ws.Range("I83").Value = 3050.1667
$ws.Range("K83").Value = 15250.8335
$ws.Range("M83").Value = -10258.8335
# row 97 (diff @ 40485)
$ws.Range("H97").Value = 1471.0834
$ws.Range("I97").Value = 1295.0625
$ws.Range("J97").Value = 1823.125
$ws.Range("K97").Value = 1295.0625
$ws.Range("L97").Value = 1823.125
$ws.Range("M97").Value = -799.0625
$ws.Range("N97").Value = -2815.125
# row 113 (diff @ 41263)
$ws.Range("H113").Value = 1951.1111
$ws.Range("I113").Value = 1352.2
$ws.Range("K113").Value = 1352.2
$ws.Range("M113").Value = 817.8
# row 122 (diff @ 41689)
$ws.Range("H122").Value = 3850.6897
$ws.Range("J122").Value = 3674.7273
$ws.Range("L122").Value = 11024.1819
$ws.Range("N122").Value = -15924.1819
# row 132 (diff @ 42179)
$ws.Range("H132").Value = 2562.0881
$ws.Range("I132").Value = 2269.9062
$ws.Range("J132").Value = 7237
$ws.Range("K132").Value = 6809.7186
$ws.Range("L132").Value = 21711
$ws.Range("M132").Value = -4279.7186
$ws.Range("N132").Value = -26771

$ws = $wb.Worksheets.Item("LTW")
# row 16 (diff @ 43440)
$ws.Range("H16").Value = 7500
$ws.Range("I16").Value = 7500
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 7500
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -7330
$ws.Range("N16").ClearContents()
# row 40 (diff @ 44646)
$ws.Range("H40").Value = 3355
$ws.Range("I40").Value = 3052.4583
$ws.Range("K40").Value = 3052.4583
$ws.Range("M40").Value = -2916.4583
# row 61 (diff @ 45675)
$ws.Range("H61").Value = 825.5
$ws.Range("I61").Value = 801
$ws.Range("K61").Value = 801
$ws.Range("M61").Value = -599
# row 113 (diff @ 48223)
$ws.Range("H113").Value = 825.5
$ws.Range("I113").Value = 801
$ws.Range("K113").Value = 801
$ws.Range("M113").Value = 1369
# row 132 (diff @ 49157)
$ws.Range("H132").Value = 2501.889
$ws.Range("I132").Value = 1814.9375
$ws.Range("J132").Value = 7997.5
$ws.Range("K132").Value = 5444.8125
$ws.Range("L132").Value = 23992.5
$ws.Range("M132").Value = -2914.8125
$ws.Range("N132").Value = -29052.5

$ws = $wb.Worksheets.Item("WVR")
# row 29 (diff @ 51064)
$ws.Range("H29").Value = 1003513
$ws.Range("I29").Value = 2003116
$ws.Range("J29").Value = 3910
$ws.Range("K29").Value = 2003116
$ws.Range("L29").Value = 3910
$ws.Range("M29").Value = -2002826
$ws.Range("N29").Value = -4490
# row 45 (diff @ 51818)
$ws.Range("H45").Value = 19313
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 19313
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 19313
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -20295
# row 122 (diff @ 55579)
$ws.Range("H122").Value = 39237.914
$ws.Range("I122").Value = 46391.45
$ws.Range("J122").Value = 4662.5
$ws.Range("K122").Value = 139174.35
$ws.Range("L122").Value = 13987.5
$ws.Range("M122").Value = -136724.35
$ws.Range("N122").Value = -18887.5
